# "Model Performances.xlsx" — document performances & fix hyperparameters
# Target sheet: "Maps 2 RGB" (internally sheet3.xml), already the active /
# tabSelected sheet in the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Mark the four V8.04.7 - V8.04.10 rows as needing revised hyperparameters ---
# (new column J, same note text on each of the 4 rows)
$ws.Range("J158:J161").Value = "REVISE hyperparameters"

# --- Newly documented model performance rows (170-197) ---

# V9.00.05
$ws.Cells.Item(170,1).Value = "V9.00.05"
$ws.Cells.Item(170,2).Value = 4.484
$ws.Cells.Item(170,3).Value = 0.5474
$ws.Cells.Item(170,4).Value = 6.5403000000000002
$ws.Cells.Item(170,5).Value = 0.69010000000000005
$ws.Cells.Item(170,6).Value = 9.4009999999999998
$ws.Cells.Item(170,7).Value = 0.73629999999999995
$ws.Cells.Item(170,8).Value = 12.702
$ws.Cells.Item(170,9).Value = 0.79690000000000005

# V9.00.06
$ws.Cells.Item(171,1).Value = "V9.00.06"
$ws.Cells.Item(171,2).Value = 4.5608000000000004
$ws.Cells.Item(171,3).Value = 0.54510000000000003
$ws.Cells.Item(171,4).Value = 6.3360000000000003
$ws.Cells.Item(171,5).Value = 0.68710000000000004
$ws.Cells.Item(171,6).Value = 9.3574000000000002
$ws.Cells.Item(171,7).Value = 0.73980000000000001
$ws.Cells.Item(171,8).Value = 13.037000000000001
$ws.Cells.Item(171,9).Value = 0.81140000000000001

# V9.00.07 - V9.00.14 : labels only (runs that didn't complete / no metrics yet)
$ws.Cells.Item(172,1).Value = "V9.00.07"
$ws.Cells.Item(173,1).Value = "V9.00.08"
$ws.Cells.Item(174,1).Value = "V9.00.09"
$ws.Cells.Item(175,1).Value = "V9.00.10"
$ws.Cells.Item(176,1).Value = "V9.00.11"
$ws.Cells.Item(177,1).Value = "V9.00.12"
$ws.Cells.Item(178,1).Value = "V9.00.13"
$ws.Cells.Item(179,1).Value = "V9.00.14"

# V9.00.15
$ws.Cells.Item(180,1).Value = "V9.00.15"
$ws.Cells.Item(180,2).Value = 4.8554000000000004
$ws.Cells.Item(180,3).Value = 0.5464
$ws.Cells.Item(180,4).Value = 6.4715999999999996
$ws.Cells.Item(180,5).Value = 0.68899999999999995
$ws.Cells.Item(180,6).Value = 8.7712000000000003
$ws.Cells.Item(180,7).Value = 0.74050000000000005
$ws.Cells.Item(180,8).Value = 13.092700000000001
$ws.Cells.Item(180,9).Value = 0.81910000000000005

# V9.00.16
$ws.Cells.Item(181,1).Value = "V9.00.16"
$ws.Cells.Item(181,2).Value = 4.6875999999999998
$ws.Cells.Item(181,3).Value = 0.53849999999999998
$ws.Cells.Item(181,4).Value = 6.1853999999999996
$ws.Cells.Item(181,5).Value = 0.68440000000000001
$ws.Cells.Item(181,6).Value = 8.6862999999999992
$ws.Cells.Item(181,7).Value = 0.73280000000000001
$ws.Cells.Item(181,8).Value = 13.4095
$ws.Cells.Item(181,9).Value = 0.82294

# V9.00.17 - V9.01.18 : labels only (remaining in-progress runs)
$ws.Cells.Item(182,1).Value = "V9.00.17"
$ws.Cells.Item(183,1).Value = "V9.00.18"
$ws.Cells.Item(184,1).Value = "V9.01.5"
$ws.Cells.Item(185,1).Value = "V9.01.6"
$ws.Cells.Item(186,1).Value = "V9.01.7"
$ws.Cells.Item(187,1).Value = "V9.01.8"
$ws.Cells.Item(188,1).Value = "V9.01.9"
$ws.Cells.Item(189,1).Value = "V9.01.10"
$ws.Cells.Item(190,1).Value = "V9.01.11"
$ws.Cells.Item(191,1).Value = "V9.01.12"
$ws.Cells.Item(192,1).Value = "V9.01.13"
$ws.Cells.Item(193,1).Value = "V9.01.14"
$ws.Cells.Item(194,1).Value = "V9.01.15"
$ws.Cells.Item(195,1).Value = "V9.01.16"
$ws.Cells.Item(196,1).Value = "V9.01.17"
$ws.Cells.Item(197,1).Value = "V9.01.18"

# --- Restore view state: selection on the working cell near the new rows ---
[void]$ws.Range("F183").Select()

# Best-effort: match the saved scroll / window geometry from the authored
# workbook (not all of these are observable through this host, but set them
# anyway in case the backing engine honours them).
try { $excel.ActiveWindow.ScrollRow = 176 } catch {}
try {
    $win = $excel.ActiveWindow
    $win.Left = 2175
    $win.Top = 1245
    $win.Width = 28800
    $win.Height = 15345
} catch {}
[void]$ws.Range("F183").Select()
